# Transmitter incidence input converted to elevation.
# The "Dynamic" sheet's column B held the transmitter incidence angle
# (Tx_th, theta). It is now re-expressed as elevation (Tx_el = 90 - theta).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dynamic")

# Rename the header for column B.
$ws.Range("B1").Value = "Tx_el (deg)"

# Find the last used row in column A (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Convert every Tx_th value in column B to elevation: el = 90 - th.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = 90 - $cell.Value2
}

# Update the active selection to M6 on the Dynamic sheet.
$ws.Activate()
$ws.Range("M6").Select()
